$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mistral")

# Extend the labeled-column formatting (bold/border/center, same as A2:A24) down to the new rows
$ws.Range("A2").Copy($ws.Range("A25:A32"))

$data = @(
    @{ Row = 2; A = 'preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse'; B = 1.30873634441749 },
    @{ Row = 3; A = 'preds_ns5_ws200_mc750_ea1.002_snks0_hopf_True_type_max_fused_opt_qcache_new_burst_lenNone_gblFalse'; B = 0 },
    @{ Row = 4; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse'; B = 0 },
    @{ Row = 5; A = 'preds_ns5_ws200_mc4000_ea1.002_snks0_hopf_True_type_max_fused_opt_qcache_new_burst_lenNone_gblFalse'; B = 1.010242314459956 },
    @{ Row = 6; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_lenNone_gblFalse'; B = 0.9920371054360548 },
    @{ Row = 7; A = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_lenNone_gblFalse'; B = 1.416201545849141 },
    @{ Row = 8; A = 'preds_ns5_ws200_mc2000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse'; B = 6.922556910629615 },
    @{ Row = 9; A = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse'; B = 0 },
    @{ Row = 10; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse'; B = 1.57627808441923 },
    @{ Row = 11; A = 'preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse'; B = 0.1295784779876364 },
    @{ Row = 12; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_lenNone_gblFalse'; B = 0.654116009700185 },
    @{ Row = 13; A = 'preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse'; B = 8.10217903827626 },
    @{ Row = 14; A = 'preds_ns5_ws200_mc2000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_sl3_lenNone_gblFalse'; B = 1.019581689577883 },
    @{ Row = 15; A = 'preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_new_lenNone_gblFalse'; B = 1.532061922169644 },
    @{ Row = 16; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_qcache_lenNone_gblFalse'; B = 1.695335097976731 },
    @{ Row = 17; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_lenNone_gblFalse'; B = 1.451674235796612 },
    @{ Row = 18; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_qcache_lenNone_gblFalse'; B = 0.8876167544943719 },
    @{ Row = 19; A = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_qcache_lenNone_gblFalse'; B = 0.7907740784882081 },
    @{ Row = 20; A = 'preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_False_type_max_fused_lenNone_gblFalse'; B = 0.6459789344721453 },
    @{ Row = 21; A = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse'; B = 4.646624851453602 },
    @{ Row = 22; A = 'preds_ns10_ws32_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse'; B = 4.511714891931707 },
    @{ Row = 23; A = 'preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_snapkv_rerun_lenNone_gblFalse'; B = 2.225652065976954 },
    @{ Row = 24; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_qcache_new_lenNone_gblFalse'; B = 1.580706269354322 },
    @{ Row = 25; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_qcache_new_lenNone_gblFalse'; B = 0.7217247351331573 },
    @{ Row = 26; A = 'preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse'; B = 0 },
    @{ Row = 27; A = 'preds_ns5_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_new_burst_lenNone_gblFalse'; B = 0 },
    @{ Row = 28; A = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_qcache_lenNone_gblFalse'; B = 1.748998216443056 },
    @{ Row = 29; A = 'preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse'; B = 3.176511576234317 },
    @{ Row = 30; A = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_new_lenNone_gblFalse'; B = 1.22951175916451 },
    @{ Row = 31; A = 'preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse'; B = 3.820181445594042 },
    @{ Row = 32; A = 'preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_sum_fused_rerun_lenNone_gblFalse'; B = 8.122498173586132 }

)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 1).Value = $row.A
    $ws.Cells.Item($row.Row, 2).Value = $row.B
}
